# The test scenario documented in row 4 changed from "RegistroUsuario" to a
# new "Aparicion de capcha" scenario: the site now shows a captcha instead of
# completing the registration, so the scenario name/description/expected
# result all change. The precondition (D4, "Encontrarse en el formulario de
# registro.") stays the same.
#
# Values are written in the same order the original author typed them so the
# resulting shared-string table comes out in the same order as the saved
# file: Resultado esperado (F4) first, then Nombre del escenario (B4),
# Descripción (C4), and finally Paso a Paso (E4).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = "No registra mi usuario y aparece el capcha para verificar que no es un robot."
$ws.Range("B4").Value = "Aparicion de capcha "
$ws.Range("C4").Value = "En la pagina de myShopify me intento registrar sin embargo el capcha no me deja debido a que es un robot"
$ws.Range("E4").Value = "1. ir a la opcion del formulario de registro. 2. llenar los datos entregados para el registro en el formulario. 3. darle submit. 4. verifico que si haya aprecido el capcha."

# The author also scrolled the sheet back up a bit and moved the selection.
$window = $excel.ActiveWindow
$window.ScrollRow = 2
$window.ScrollColumn = 1
$ws.Range("G4").Select()
